$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33 (ALC)
$ws.Range("H33").Value = 19608382
$ws.Range("I33").Value = 30303652
$ws.Range("J33").Value = 388.66666
$ws.Range("K33").Value = 30303652
$ws.Range("L33").Value = 388.66666
$ws.Range("M33").Value = -30303423
$ws.Range("N33").Value = -846.66666

# Row 42 (ALC)
$ws.Range("H42").Value = 453.18182
$ws.Range("I42").Value = 96
$ws.Range("J42").Value = 587.125
$ws.Range("K42").Value = 288
$ws.Range("L42").Value = 1761.375
$ws.Range("M42").Value = -58
$ws.Range("N42").Value = -2221.375

# Row 62 (ALC)
$ws.Range("H62").Value = 2134.6667
$ws.Range("I62").Value = 2362
$ws.Range("J62").Value = 1850.5
$ws.Range("K62").Value = 2362
$ws.Range("L62").Value = 1850.5
$ws.Range("M62").Value = -1738
$ws.Range("N62").Value = -3098.5

# Row 65 (ALC)
$ws.Range("H65").Value = 2134.6667
$ws.Range("I65").Value = 2362
$ws.Range("J65").Value = 1850.5
$ws.Range("K65").Value = 11810
$ws.Range("L65").Value = 9252.5
$ws.Range("M65").Value = -8690
$ws.Range("N65").Value = -15492.5

# Row 106 (ALC)
$ws.Range("H106").Value = 20877108
$ws.Range("I106").Value = 51879
$ws.Range("J106").Value = 125003250
$ws.Range("K106").Value = 51879
$ws.Range("L106").Value = 125003250
$ws.Range("M106").Value = -51248
$ws.Range("N106").Value = -125004512

# Row 132 (ALC)
$ws.Range("H132").Value = 1809893
$ws.Range("I132").Value = 1906295.2
$ws.Range("J132").Value = 2350
$ws.Range("K132").Value = 5718885.6
$ws.Range("L132").Value = 7050
$ws.Range("M132").Value = -5716355.6
$ws.Range("N132").Value = -12110

# Row 134 (ALC)
$ws.Range("H134").Value = 38000
$ws.Range("J134").Value = 38000
$ws.Range("L134").Value = 38000
$ws.Range("N134").Value = -48140

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 33500.605
$ws.Range("I32").Value = 37244.8
$ws.Range("J32").Value = 28460.346
$ws.Range("K32").Value = 37244.8
$ws.Range("L32").Value = 28460.346
$ws.Range("M32").Value = -36957.8
$ws.Range("N32").Value = -29034.346

# Row 45 (ARM)
$ws.Range("H45").Value = 973.3333
$ws.Range("I45").Value = 973.3333
$ws.Range("K45").Value = 973.3333
$ws.Range("M45").Value = -596.3333

$ws = $wb.Worksheets.Item("BSM")
# Row 134 (BSM)
$ws.Range("H134").Value = 4907
$ws.Range("I134").Value = 6625.2593
$ws.Range("J134").Value = 2697.8096
$ws.Range("K134").Value = 19875.7779
$ws.Range("L134").Value = 8093.4288
$ws.Range("M134").Value = -17340.7779
$ws.Range("N134").Value = -13163.4288

$ws = $wb.Worksheets.Item("CRP")
# Row 107 (CRP)
$ws.Range("H107").Value = 714.4167
$ws.Range("I107").Value = 731.6111
$ws.Range("K107").Value = 731.6111
$ws.Range("M107").Value = 1188.3889

# Row 134 (CRP)
$ws.Range("H134").Value = 5056.48
$ws.Range("I134").Value = 5708.5713
$ws.Range("J134").Value = 1633
$ws.Range("K134").Value = 17125.7139
$ws.Range("L134").Value = 4899
$ws.Range("M134").Value = -14590.7139
$ws.Range("N134").Value = -9969

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (CUL)
$ws.Range("H5").Value = 989.62164
$ws.Range("I5").Value = 294.92856
$ws.Range("J5").Value = 3150.889
$ws.Range("K5").Value = 884.78568
$ws.Range("L5").Value = 9452.667000000001
$ws.Range("M5").Value = -772.78568
$ws.Range("N5").Value = -9676.667000000001

# Row 80 (CUL)
$ws.Range("H80").Value = 5000
$ws.Range("I80").Value = 1000
$ws.Range("J80").Value = 6000
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 18000
$ws.Range("M80").Value = -2064
$ws.Range("N80").Value = -19872

# Row 83 (CUL)
$ws.Range("H83").Value = 5000
$ws.Range("I83").Value = 1000
$ws.Range("J83").Value = 6000
$ws.Range("K83").Value = 9000
$ws.Range("L83").Value = 54000
$ws.Range("M83").Value = -4320
$ws.Range("N83").Value = -63360

# Row 86 (CUL)
$ws.Range("H86").Value = 2217.6667
$ws.Range("I86").Value = 1910.6875
$ws.Range("J86").Value = 3200
$ws.Range("K86").Value = 5732.0625
$ws.Range("L86").Value = 9600
$ws.Range("M86").Value = -4546.0625
$ws.Range("N86").Value = -11972

# Row 89 (CUL)
$ws.Range("H89").Value = 2217.6667
$ws.Range("I89").Value = 1910.6875
$ws.Range("J89").Value = 3200
$ws.Range("K89").Value = 17196.1875
$ws.Range("L89").Value = 28800
$ws.Range("M89").Value = -11268.1875
$ws.Range("N89").Value = -40656

# Row 131 (CUL)
$ws.Range("H131").Value = 1607559.6
$ws.Range("J131").Value = 1725284.6
$ws.Range("L131").Value = 5175853.800000001
$ws.Range("N131").Value = -5185933.800000001

# Row 135 (CUL)
$ws.Range("H135").Value = 989.62164
$ws.Range("I135").Value = 294.92856
$ws.Range("J135").Value = 3150.889
$ws.Range("K135").Value = 2654.35704
$ws.Range("L135").Value = 28358.001
$ws.Range("M135").Value = -119.3570399999999
$ws.Range("N135").Value = -33428.001

# Row 137 (CUL)
$ws.Range("H137").Value = 55577684
$ws.Range("I137").Value = 41681092
$ws.Range("J137").Value = 62117256
$ws.Range("K137").Value = 125043276
$ws.Range("L137").Value = 186351768
$ws.Range("M137").Value = -125038176
$ws.Range("N137").Value = -186361968

$ws = $wb.Worksheets.Item("GSM")
# Row 113 (GSM)
$ws.Range("H113").Value = 25001054
$ws.Range("I113").Value = 35715148
$ws.Range("J113").Value = 1504.3334
$ws.Range("K113").Value = 35715148
$ws.Range("L113").Value = 1504.3334
$ws.Range("M113").Value = -35712978
$ws.Range("N113").Value = -5844.3334

# Row 132 (GSM)
$ws.Range("H132").Value = 7015.2085
$ws.Range("I132").Value = 11100.333
$ws.Range("J132").Value = 2930.0833
$ws.Range("K132").Value = 33300.999
$ws.Range("L132").Value = 8790.249899999999
$ws.Range("M132").Value = -30770.999
$ws.Range("N132").Value = -13850.2499

$ws = $wb.Worksheets.Item("LTW")
# Row 68 (LTW)
$ws.Range("H68").Value = 125000000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 125000000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 125000000
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -125001498

# Row 71 (LTW)
$ws.Range("H71").Value = 125000000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 125000000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 625000000
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -625007488

# Row 132 (LTW)
$ws.Range("H132").Value = 7421.8857
$ws.Range("I132").Value = 10805.9
$ws.Range("J132").Value = 2909.8667
$ws.Range("K132").Value = 32417.7
$ws.Range("L132").Value = 8729.6001
$ws.Range("M132").Value = -29887.7
$ws.Range("N132").Value = -13789.6001

# Row 136 (LTW)
$ws.Range("H136").Value = 5188.0835
$ws.Range("I136").Value = 6071.769
$ws.Range("J136").Value = 2890.5
$ws.Range("K136").Value = 18215.307
$ws.Range("L136").Value = 8671.5
$ws.Range("M136").Value = -15665.307
$ws.Range("N136").Value = -13771.5

$ws = $wb.Worksheets.Item("WVR")
# Row 136 (WVR)
$ws.Range("H136").Value = 1701.7556
$ws.Range("I136").Value = 1669.7567
$ws.Range("K136").Value = 5009.2701
$ws.Range("M136").Value = -2459.2701
